$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7476426666666667
$ws.Range("H2").Value = 2.242928
$ws.Range("I2").Value = 0.001581772089386036
$ws.Range("J2").Value = 0.001581772089386036
$ws.Range("M2").Value = 9.423852333333334
$ws.Range("N2").Value = 28.271557
$ws.Range("O2").Value = 0.06654336290212845
$ws.Range("P2").Value = 0.06654336290212845
$ws.Range("Q2").Value = 7.045674088766223
$ws.Range("R2").Value = 63.411066798896
$ws.Range("S2").Value = 0.0001052564341724729
$ws.Range("T2").Value = 0.000105256434172473
$ws.Range("G3").Value = 0.7476426666666667
$ws.Range("H3").Value = 2.242928
$ws.Range("I3").Value = 0.001581772089386036
$ws.Range("J3").Value = 0.001581772089386036
$ws.Range("O3").Value = 0.3572423751649123
$ws.Range("P3").Value = 0.3572423751649123
$ws.Range("Q3").Value = 37.82515996089245
$ws.Range("R3").Value = 340.426439648032
$ws.Range("S3").Value = 0.0005650760181818334
$ws.Range("T3").Value = 0.0005650760181818335
$ws.Range("G4").Value = 0.7476426666666667
$ws.Range("H4").Value = 2.242928
$ws.Range("I4").Value = 0.001581772089386036
$ws.Range("J4").Value = 0.001581772089386036
$ws.Range("M4").Value = 26.84076266666667
$ws.Range("N4").Value = 80.522288
$ws.Range("O4").Value = 0.1895270158659356
$ws.Range("P4").Value = 0.1895270158659356
$ws.Range("Q4").Value = 20.06729937547378
$ws.Range("R4").Value = 180.605694379264
$ws.Range("S4").Value = 0.0002997885438813613
$ws.Range("T4").Value = 0.0002997885438813613
$ws.Range("G5").Value = 0.7476426666666667
$ws.Range("H5").Value = 2.242928
$ws.Range("I5").Value = 0.001581772089386036
$ws.Range("J5").Value = 0.001581772089386036
$ws.Range("M5").Value = 54.762539
$ws.Range("N5").Value = 164.287617
$ws.Range("O5").Value = 0.3866872460670236
$ws.Range("P5").Value = 0.3866872460670236
$ws.Range("Q5").Value = 40.94281069139733
$ws.Range("R5").Value = 368.485296222576
$ws.Range("S5").Value = 0.0006116510931503681
$ws.Range("T5").Value = 0.0006116510931503682
$ws.Range("I6").Value = 0.02590993131491687
$ws.Range("J6").Value = 0.02590993131491688
$ws.Range("M6").Value = 9.423852333333334
$ws.Range("N6").Value = 28.271557
$ws.Range("O6").Value = 0.06654336290212845
$ws.Range("P6").Value = 0.06654336290212845
$ws.Range("Q6").Value = 115.4103887229925
$ws.Range("R6").Value = 1038.693498506932
$ws.Range("S6").Value = 0.001724133962257736
$ws.Range("T6").Value = 0.001724133962257736
$ws.Range("I7").Value = 0.02590993131491687
$ws.Range("J7").Value = 0.02590993131491688
$ws.Range("O7").Value = 0.3572423751649123
$ws.Range("P7").Value = 0.3572423751649123
$ws.Range("Q7").Value = 619.5881841251049
$ws.Range("R7").Value = 5576.293657125944
$ws.Range("S7").Value = 0.009256125403300643
$ws.Range("T7").Value = 0.009256125403300645
$ws.Range("I8").Value = 0.02590993131491687
$ws.Range("J8").Value = 0.02590993131491688
$ws.Range("M8").Value = 26.84076266666667
$ws.Range("N8").Value = 80.522288
$ws.Range("O8").Value = 0.1895270158659356
$ws.Range("P8").Value = 0.1895270158659356
$ws.Range("Q8").Value = 328.7087640395875
$ws.Range("R8").Value = 2958.378876356288
$ws.Range("S8").Value = 0.004910631963407551
$ws.Range("T8").Value = 0.004910631963407553
$ws.Range("I9").Value = 0.02590993131491687
$ws.Range("J9").Value = 0.02590993131491688
$ws.Range("M9").Value = 54.762539
$ws.Range("N9").Value = 164.287617
$ws.Range("O9").Value = 0.3866872460670236
$ws.Range("P9").Value = 0.3866872460670236
$ws.Range("Q9").Value = 670.6562974350546
$ws.Range("R9").Value = 6035.906676915492
$ws.Range("S9").Value = 0.01001903998595094
$ws.Range("T9").Value = 0.01001903998595094
$ws.Range("G10").Value = 18.93023433333333
$ws.Range("H10").Value = 56.79070299999999
$ws.Range("I10").Value = 0.04005030430848061
$ws.Range("J10").Value = 0.04005030430848062
$ws.Range("M10").Value = 9.423852333333334
$ws.Range("N10").Value = 28.271557
$ws.Range("O10").Value = 0.06654336290212845
$ws.Range("P10").Value = 0.06654336290212845
$ws.Range("Q10").Value = 178.3957329927301
$ws.Range("R10").Value = 1605.561596934571
$ws.Range("S10").Value = 0.002665081933939904
$ws.Range("T10").Value = 0.002665081933939904
$ws.Range("G11").Value = 18.93023433333333
$ws.Range("H11").Value = 56.79070299999999
$ws.Range("I11").Value = 0.04005030430848061
$ws.Range("J11").Value = 0.04005030430848062
$ws.Range("O11").Value = 0.3572423751649123
$ws.Range("P11").Value = 0.3572423751649123
$ws.Range("Q11").Value = 957.7291046643202
$ws.Range("R11").Value = 8619.56194197888
$ws.Range("S11").Value = 0.01430766583723913
$ws.Range("T11").Value = 0.01430766583723914
$ws.Range("G12").Value = 18.93023433333333
$ws.Range("H12").Value = 56.79070299999999
$ws.Range("I12").Value = 0.04005030430848061
$ws.Range("J12").Value = 0.04005030430848062
$ws.Range("M12").Value = 26.84076266666667
$ws.Range("N12").Value = 80.522288
$ws.Range("O12").Value = 0.1895270158659356
$ws.Range("P12").Value = 0.1895270158659356
$ws.Range("Q12").Value = 508.1019269653848
$ws.Range("R12").Value = 4572.917342688464
$ws.Range("S12").Value = 0.007590614660108953
$ws.Range("T12").Value = 0.007590614660108954
$ws.Range("G13").Value = 18.93023433333333
$ws.Range("H13").Value = 56.79070299999999
$ws.Range("I13").Value = 0.04005030430848061
$ws.Range("J13").Value = 0.04005030430848062
$ws.Range("M13").Value = 54.762539
$ws.Range("N13").Value = 164.287617
$ws.Range("O13").Value = 0.3866872460670236
$ws.Range("P13").Value = 0.3866872460670236
$ws.Range("Q13").Value = 1036.667695958306
$ws.Range("R13").Value = 9330.009263624748
$ws.Range("S13").Value = 0.01548694187719262
$ws.Range("T13").Value = 0.01548694187719262
$ws.Range("G14").Value = 440.7369333333333
$ws.Range("H14").Value = 1322.2108
$ws.Range("I14").Value = 0.9324579922872165
$ws.Range("J14").Value = 0.9324579922872166
$ws.Range("M14").Value = 9.423852333333334
$ws.Range("N14").Value = 28.271557
$ws.Range("O14").Value = 0.06654336290212845
$ws.Range("P14").Value = 0.06654336290212845
$ws.Range("Q14").Value = 4153.439777579511
$ws.Range("R14").Value = 37380.9579982156
$ws.Range("S14").Value = 0.06204889057175834
$ws.Range("T14").Value = 0.06204889057175834
$ws.Range("G15").Value = 440.7369333333333
$ws.Range("H15").Value = 1322.2108
$ws.Range("I15").Value = 0.9324579922872165
$ws.Range("J15").Value = 0.9324579922872166
$ws.Range("O15").Value = 0.3572423751649123
$ws.Range("P15").Value = 0.3572423751649123
$ws.Range("Q15").Value = 22298.01180065502
$ws.Range("R15").Value = 200682.1062058952
$ws.Range("S15").Value = 0.3331135079061907
$ws.Range("T15").Value = 0.3331135079061908
$ws.Range("G16").Value = 440.7369333333333
$ws.Range("H16").Value = 1322.2108
$ws.Range("I16").Value = 0.9324579922872165
$ws.Range("J16").Value = 0.9324579922872166
$ws.Range("M16").Value = 26.84076266666667
$ws.Range("N16").Value = 80.522288
$ws.Range("O16").Value = 0.1895270158659356
$ws.Range("P16").Value = 0.1895270158659356
$ws.Range("Q16").Value = 11829.71542603449
$ws.Range("R16").Value = 106467.4388343104
$ws.Range("S16").Value = 0.1767259806985377
$ws.Range("T16").Value = 0.1767259806985378
$ws.Range("G17").Value = 440.7369333333333
$ws.Range("H17").Value = 1322.2108
$ws.Range("I17").Value = 0.9324579922872165
$ws.Range("J17").Value = 0.9324579922872166
$ws.Range("M17").Value = 54.762539
$ws.Range("N17").Value = 164.287617
$ws.Range("O17").Value = 0.3866872460670236
$ws.Range("P17").Value = 0.3866872460670236
$ws.Range("Q17").Value = 24135.87350040706
$ws.Range("R17").Value = 217222.8615036636
$ws.Range("S17").Value = 0.3605696131107297
$ws.Range("T17").Value = 0.3605696131107298
